$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I8 originally held " 6.19"; the commit doubles it to " 6.19 + 6.19" (6.19 -> 12.38).
# Everything else (J8, J13, I16, H18) is recomputed automatically by the
# engine's dependency chain once I8 changes.
$ws.Range("I8").Formula = "= 6.19 + 6.19"
